$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings) ---
# A8: "Volume 30   Number  37" -> "...38"
$ws.Range("A8").Characters(21, 2).Text = "38"

# C9: "Report Covering the Week  9/11/2023  Through  9/17/2023"
#     -> "...9/18/2023  Through  9/24/2023"
$ws.Range("C9").Characters(27, 9).Text = "9/18/2023"
$ws.Range("C9").Characters(47, 9).Text = "9/24/2023"

# --- Crime statistics table updates (rows 14-29) ---
$ws.Range("M14").Value = -81.818181818181
$ws.Range("F15").Value = 4
$ws.Range("H15").Value = 300
$ws.Range("I15").Value = 15
$ws.Range("K15").Value = 66.666666666666
$ws.Range("L15").Value = -54.545454545454
$ws.Range("M15").Value = 7.142857142857
$ws.Range("N15").Value = -54.545454545454
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("D16").NumberFormat = "#,##0"
$ws.Range("E16").Value = -28.571428571428
$ws.Range("E16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -6.25
$ws.Range("I16").Value = 126
$ws.Range("J16").Value = 140
$ws.Range("K16").Value = -10
$ws.Range("L16").Value = 18.867924528301
$ws.Range("M16").Value = -47.058823529411
$ws.Range("N16").Value = -83.244680851063
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 38
$ws.Range("G17").Value = 28
$ws.Range("H17").Value = 35.714285714285
$ws.Range("I17").Value = 333
$ws.Range("J17").Value = 324
$ws.Range("K17").Value = 2.777777777777
$ws.Range("L17").Value = 10.264900662251
$ws.Range("M17").Value = 57.075471698113
$ws.Range("N17").Value = 10.264900662251
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 20
$ws.Range("F18").Value = 24
$ws.Range("G18").Value = 28
$ws.Range("H18").Value = -14.285714285714
$ws.Range("I18").Value = 168
$ws.Range("J18").Value = 145
$ws.Range("K18").Value = 15.862068965517
$ws.Range("L18").Value = 57.009345794392
$ws.Range("M18").Value = -36.842105263157
$ws.Range("N18").Value = -85.762711864406
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -15.384615384615
$ws.Range("F19").Value = 52
$ws.Range("G19").Value = 41
$ws.Range("H19").Value = 26.829268292682
$ws.Range("I19").Value = 509
$ws.Range("J19").Value = 469
$ws.Range("K19").Value = 8.528784648187
$ws.Range("L19").Value = 58.07453416149
$ws.Range("M19").Value = 59.0625
$ws.Range("N19").Value = 14.639639639639
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 39
$ws.Range("G20").Value = 26
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 257
$ws.Range("J20").Value = 182
$ws.Range("K20").Value = 41.208791208791
$ws.Range("L20").Value = 142.452830188679
$ws.Range("M20").Value = -9.187279151943
$ws.Range("N20").Value = -89.641273679967
$ws.Range("C21").Value = 41
$ws.Range("D21").Value = 43
$ws.Range("E21").Value = -4.651162790697
$ws.Range("F21").Value = 172
$ws.Range("G21").Value = 141
$ws.Range("H21").Value = 21.985815602836
$ws.Range("I21").Value = 1410
$ws.Range("J21").Value = 1273
$ws.Range("K21").Value = 10.761979575805
$ws.Range("L21").Value = 43.877551020408
$ws.Range("M21").Value = 4.910714285714
$ws.Range("N21").Value = -72.915866308106
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = -22.58064516129
$ws.Range("F24").Value = 107
$ws.Range("G24").Value = 114
$ws.Range("H24").Value = -6.140350877192
$ws.Range("I24").Value = 938
$ws.Range("J24").Value = 1124
$ws.Range("K24").Value = -16.548042704626
$ws.Range("L24").Value = 34.383954154727
$ws.Range("M24").Value = 50.08
$ws.Range("C25").Value = 26
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 550
$ws.Range("F25").Value = 71
$ws.Range("G25").Value = 50
$ws.Range("H25").Value = 42
$ws.Range("I25").Value = 494
$ws.Range("J25").Value = 480
$ws.Range("K25").Value = 2.916666666666
$ws.Range("L25").Value = 40.74074074074
$ws.Range("M25").Value = 6.465517241379
$ws.Range("C26").Value = 1
$ws.Range("I26").Value = 27
$ws.Range("K26").Value = 28.571428571428
$ws.Range("L26").Value = -35.714285714285
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -16.666666666666
$ws.Range("I27").Value = 37
$ws.Range("J27").Value = 43
$ws.Range("K27").Value = -13.953488372093
$ws.Range("L27").Value = -5.128205128205
$ws.Range("D28").Value = 1
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("E28").Value = -100
$ws.Range("E28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G28").Value = 2
$ws.Range("J28").Value = 20
$ws.Range("K28").Value = -65
$ws.Range("M28").Value = -78.125
$ws.Range("D29").Value = 1
$ws.Range("D29").NumberFormat = "#,##0"
$ws.Range("E29").Value = -100
$ws.Range("E29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G29").Value = 2
$ws.Range("J29").Value = 16
$ws.Range("K29").Value = -62.5
$ws.Range("M29").Value = -73.91304347826
